# Updated cryptos list on Wed May 17 05:58:05 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.030.22"
$ws.Range("E2").Value = "  -0.65%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.821.81"
$ws.Range("E3").Value = "  +0.03%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.54%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.81%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.47%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4547"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.41%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3710"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.14%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07285"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.80%  "

# Row 10 - Polygon
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8576"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.40%  "

# Row 11 - Solana
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.80"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.84%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.821.66"
$ws.Range("E12").Value = "  -0.90%  "

# Row 13 - Chainlink
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.654"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.23%  "

# Row 14 - Litecoin
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.16%  "

# Row 15 - Polkadot -> TRON (rank swap)
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07101"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.10%  "

# Row 16 - TRON -> Polkadot (rank swap)
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.327"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.35%  "

# Row 17 - BinanceUSD
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.63%  "

# Row 18 - ShibaInu
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008810"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.53%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  -0.46%  "

# Row 20 - Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.56%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "27.070.77"
$ws.Range("E21").Value = "  -0.61%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +0.66%  "

# Row 23 - Cosmos
$ws.Range("E23").Value = "  +0.66%  "

# Row 24 - Toncoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.990"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.69%  "

# Row 25 - Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.31%  "

# Row 26 - LidoDAOToken
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.214"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.77%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  +0.40%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.246"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.42%  "

# Row 29 - BitcoinCash
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.49%  "

# Row 30 - Stellar
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08887"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.16%  "

# Row 31 - ARBITRUM
$ws.Range("E31").Value = "  -0.34%  "

# Row 32 - ImmutableX
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7524"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.70%  "

# Row 33 - HuobiToken
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.941"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.58%  "

# Row 34 - Filecoin
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.458"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "

# Row 35 - Frax
$ws.Range("E35").Value = "  -0.53%  "

# Row 36 - TrustWalletToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.098"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.54%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  -0.15%  "

# Row 38 - Hedera
$ws.Range("E38").Value = "  -0.31%  "

# Row 39 - TheSandbox
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5317"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.50%  "

# Row 40 - FraxShare
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.201"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.74%  "

# Row 41 - MXToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.876"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.31%  "

# Row 42 - Algorand
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1713"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.62%  "

# Row 43 - Decentraland
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5202"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.61%  "

# Row 44 - Aptos
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.546"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.69%  "

# Row 45 - EnergySwap
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.67"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.83%  "

# Row 46 - RenderToken
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.969"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.20%  "

# Row 47 - Quant
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.03%  "

# Row 48 - NEARProtocol -> PaxDollar (rank swap)
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.53%  "

# Row 49 - PaxDollar -> NEARProtocol (rank swap)
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.670"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.67%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  +0.13%  "

# Row 51 - Aave
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.44%  "
